$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price column (D) cells whose new values look numeric,
# so Excel stores them as text (matching the original inlineStr string type)
# instead of silently converting them to numbers.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Apply the updated cell values
$ws.Range('D2').Value = '63.558.13'
$ws.Range('E2').Value = '  +6.03%  '
$ws.Range('D3').Value = '3.400.22'
$ws.Range('E3').Value = '  +6.65%  '
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').Value = '577.19'
$ws.Range('E5').Value = '  +7.56%  '
$ws.Range('D6').Value = '155.10'
$ws.Range('E6').Value = '  +7.04%  '
$ws.Range('D7').Value = '1.00'
$ws.Range('E7').Value = '  +0.03%  '
$ws.Range('D8').Value = '3.418.00'
$ws.Range('E8').Value = '  +7.04%  '
$ws.Range('D9').Value = '0.536'
$ws.Range('E9').Value = '  +0.64%  '
$ws.Range('E10').Value = '  +2.99%  '
$ws.Range('E11').Value = '  +7.54%  '
$ws.Range('D12').Value = '0.435'
$ws.Range('E12').Value = '  +0.41%  '
$ws.Range('D13').Value = '3.987.46'
$ws.Range('E14').Value = '  +0.24%  '
$ws.Range('E15').Value = '  +7.68%  '
$ws.Range('D16').Value = '27.10'
$ws.Range('E16').Value = '  +5.44%  '
$ws.Range('D17').Value = '63.636.86'
$ws.Range('E17').Value = '  +6.18%  '
$ws.Range('D18').Value = '3.394.47'
$ws.Range('E18').Value = '  +6.39%  '
$ws.Range('D19').Value = '6.38'
$ws.Range('E19').Value = '  +2.38%  '
$ws.Range('D20').Value = '14.10'
$ws.Range('E20').Value = '  +6.55%  '
$ws.Range('D21').Value = '8.42'
$ws.Range('E21').Value = '  +2.88%  '
$ws.Range('D22').Value = '388.34'
$ws.Range('E22').Value = '  +5.24%  '
$ws.Range('D23').Value = '0.996'
$ws.Range('E23').Value = '  -0.36%  '
$ws.Range('D24').Value = '0.536'
$ws.Range('E24').Value = '  +2.54%  '
$ws.Range('E25').Value = '  +2.37%  '
$ws.Range('B26').Value = 'PEPE'
$ws.Range('C26').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D26').Value = '0.0000106'
$ws.Range('E26').Value = '  +21.53%  '
$ws.Range('B27').Value = 'InternetComputer(DFINITY)'
$ws.Range('C27').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D27').Value = '9.52'
$ws.Range('E27').Value = '  +10.87%  '
$ws.Range('E28').Value = '  +7.18%  '
$ws.Range('E29').Value = '  +0.74%  '
$ws.Range('E30').Value = '  +8.09%  '
$ws.Range('D31').Value = '6.58'
$ws.Range('E31').Value = '  +8.25%  '
$ws.Range('D32').Value = '5.73'
$ws.Range('E32').Value = '  +8.56%  '
$ws.Range('E33').Value = '  +13.47%  '
$ws.Range('D34').Value = '23.25'
$ws.Range('E34').Value = '  +3.52%  '
$ws.Range('E35').Value = '  -0.08%  '
$ws.Range('D36').Value = '6.72'
$ws.Range('E36').Value = '  +2.29%  '
$ws.Range('E37').Value = '  +9.35%  '
$ws.Range('D38').Value = '158.25'
$ws.Range('E38').Value = '  +0.03%  '
$ws.Range('D39').Value = '1.88'
$ws.Range('E39').Value = '  +11.61%  '
$ws.Range('B40').Value = 'Hedera'
$ws.Range('C40').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D40').Value = '0.0774'
$ws.Range('E40').Value = '  +9.18%  '
$ws.Range('B41').Value = 'EnergySwap'
$ws.Range('C41').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D41').Value = '27.60'
$ws.Range('E41').Value = '  +4.93%  '
$ws.Range('D42').Value = '2.910.51'
$ws.Range('E42').Value = '  +4.51%  '
$ws.Range('D43').Value = '0.0320'
$ws.Range('E43').Value = '  +3.74%  '
$ws.Range('E44').Value = '  +6.72%  '
$ws.Range('D46').Value = '4.32'
$ws.Range('E46').Value = '  +2.73%  '
$ws.Range('D47').Value = '1.08'
$ws.Range('E47').Value = '  +9.97%  '
$ws.Range('D48').Value = '22.55'
$ws.Range('E48').Value = '  +9.75%  '
$ws.Range('D49').Value = '3.450.13'
$ws.Range('E49').Value = '  +6.83%  '
$ws.Range('B50').Value = 'Bittensor'
$ws.Range('C50').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D50').Value = '295.82'
$ws.Range('E50').Value = '  +12.93%  '
$ws.Range('B51').Value = 'Cosmos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D51').Value = '6.34'
$ws.Range('E51').Value = '  +3.24%  '
